# Re-sort the data rows (2..54) of the "Avverkningsanmälningar" sheet into the
# new order supplied by the upstream data refresh, and bump the "Förändrad"
# (column C) timestamp from 2026-02-26 (serial 46079) to 2026-02-28 (serial 46081)
# for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 54

# Mapping: new row number -> old row number (where its data used to live).
$map = @{
    2=2; 3=3; 4=4; 5=6; 6=8; 7=5; 8=7; 9=9; 10=12; 11=10; 12=11; 13=13; 14=14;
    15=15; 16=19; 17=20; 18=24; 19=25; 20=26; 21=27; 22=28; 23=29; 24=31; 25=32;
    26=33; 27=16; 28=17; 29=18; 30=22; 31=44; 32=47; 33=34; 34=35; 35=43; 36=36;
    37=45; 38=37; 39=38; 40=40; 41=41; 42=42; 43=49; 44=50; 45=46; 46=48; 47=51;
    48=54; 49=52; 50=53; 51=21; 52=23; 53=30; 54=39
}

# Snapshot every data row's plain values (columns A-R) and formulas (columns S-Y)
# before anything is written back, so the permutation can be applied safely even
# though source and destination ranges overlap.
#
# NOTE: do not use "A$r:R$r" style interpolation here - in this PowerShell the
# "$r:" sequence is parsed as a scoped-variable reference (like $r:Y) rather
# than the variable followed by a literal colon, silently truncating the
# address. Build the address with string concatenation instead.
$valSnapshot = @{}
$formulaSnapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $valAddr = "A" + $r + ":R" + $r
    $formulaAddr = "S" + $r + ":Y" + $r
    $valSnapshot[$r] = $ws.Range($valAddr).Value()
    $formulaSnapshot[$r] = $ws.Range($formulaAddr).Formula()
}

for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $map[$newRow]

    $valAddr = "A" + $newRow + ":R" + $newRow
    $formulaAddr = "S" + $newRow + ":Y" + $newRow

    $ws.Range($valAddr).Value = $valSnapshot[$oldRow]
    $ws.Range($formulaAddr).Formula = $formulaSnapshot[$oldRow]

    # "Förändrad" column always becomes 2026-02-28 (serial 46081) for every row.
    $ws.Cells.Item($newRow, 3).Value = 46081
}
